$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row just above the current row 357. This shifts the
# existing rows 357..451 down to 358..452 (Excel's default Insert behaviour:
# cells shift down, formatting is inherited from the row above), matching
# the OOXML diff where every row from 357 onward is replaced by the row that
# used to precede it, and a brand-new row of weekly data lands at 357 while
# the workbook's used range grows from A1:R451 to A1:R452.
$ws.Rows(357).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(357, 1).Value = 8
$ws.Cells.Item(357, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(357, 3).Value = "Coquimbo"
$ws.Cells.Item(357, 4).Value = 45135
$ws.Cells.Item(357, 5).Value = 4
$ws.Cells.Item(357, 6).Value = 100112012
$ws.Cells.Item(357, 7).Value = "Espinaca"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 1160
$ws.Cells.Item(357, 11).Value = 500
$ws.Cells.Item(357, 12).Value = 600
$ws.Cells.Item(357, 13).Value = 550
$ws.Cells.Item(357, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(357, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(357, 16).Value = 1100
$ws.Cells.Item(357, 17).Value = 0.5
$ws.Cells.Item(357, 18).Value = "Hortaliza"
